$d = $word.ActiveDocument

# --- Part 1: split the last run of the "Regression testing" paragraph ---
# Original combined text "...на бъгс" (two runs: "...бъг" + "с") becomes
# "...на " (Cyrillic run) + "bug" (new run, en-US).
$findRng = $d.Content
$found = $findRng.Find.Execute("бъгс", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find target text 'бъгс' in document"
}
$bugTarget = $d.Range($findRng.Start, $findRng.End)
$bugXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>bug</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$bugTarget.InsertXML($bugXml, "Replace")

# --- Part 2: append 9 new list paragraphs after the "Regression testing" paragraph ---
$lastPara = $d.Paragraphs.Last
$insertionPoint = $d.Range($lastPara.Range.End, $lastPara.Range.End)
$parasXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Bug tracking systems -  github, JIRA</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Test driven development -  </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>първо пишем тестовете и после кода</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Unit Testing (white testing) – </w:t></w:r><w:r><w:t>гледаме кода като тествам, пишат се от програмистите</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> Integration test – </w:t></w:r><w:r><w:t>тестваме взаимодействието между компонентите, два или повече</w:t></w:r><w:r><w:t xml:space="preserve">(от </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">QA </w:t></w:r><w:r><w:t xml:space="preserve">или </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> Developer</w:t></w:r><w:r><w:t>)</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">System testing – </w:t></w:r><w:r><w:t xml:space="preserve">истинско проверяване на проекта(цялата система), може да тестваме със </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Selenium</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Test triangle </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Система за тестване - </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Test Rail</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Test Link, MS Test Manager,</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> Framework </w:t></w:r><w:r><w:t xml:space="preserve">за тетсване </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>C# - NUnit, Selenium</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertionPoint.InsertXML($parasXml, "After")
